$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.987.52'
$ws.Range("E2").Value = '  -1.54%  '
$ws.Range("D3").Value = '1.838.71'
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.57%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.69'
$ws.Range("E5").Value = '  -0.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4612'
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3675'
$ws.Range("E8").Value = '  -1.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07225'
$ws.Range("E9").Value = '  -2.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8808'
$ws.Range("E10").Value = '  -1.22%  '
$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").Value = '1.988.29'
$ws.Range("E11").Value = '  +8.55%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07826'
$ws.Range("E12").Value = '  -1.79%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.73'
$ws.Range("E13").Value = '  -2.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.355'
$ws.Range("E14").Value = '  -1.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.455'
$ws.Range("E15").Value = '  -2.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.25'
$ws.Range("E16").Value = '  -2.84%  '
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008781'
$ws.Range("E18").Value = '  -2.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.006'
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").Value = '27.029.64'
$ws.Range("E20").Value = '  -1.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.54'
$ws.Range("E21").Value = '  -2.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.004'
$ws.Range("E22").Value = '  -3.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.46'
$ws.Range("E23").Value = '  -1.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.002'
$ws.Range("E24").Value = '  +7.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.49'
$ws.Range("E25").Value = '  -1.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.30'
$ws.Range("E26").Value = '  -1.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.001'
$ws.Range("E27").Value = '  -4.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.01'
$ws.Range("E28").Value = '  -2.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.947'
$ws.Range("E29").Value = '  -4.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08841'
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.144'
$ws.Range("E31").Value = '  +5.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7637'
$ws.Range("E32").Value = '  +0.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.470'
$ws.Range("E33").Value = '  -1.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.136'
$ws.Range("E34").Value = '  -2.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.643'
$ws.Range("E35").Value = '  +1.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.098'
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01928'
$ws.Range("E37").Value = '  -1.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05163'
$ws.Range("E38").Value = '  -2.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.934'
$ws.Range("E39").Value = '  -2.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.977'
$ws.Range("E40").Value = '  -3.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4988'
$ws.Range("E41").Value = '  -4.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1602'
$ws.Range("E42").Value = '  -2.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.353'
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4680'
$ws.Range("E44").Value = '  -4.91%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.24'
$ws.Range("E45").Value = '  -0.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.006'
$ws.Range("E46").Value = '  +0.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.65'
$ws.Range("E47").Value = '  -0.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.614'
$ws.Range("E48").Value = '  -2.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06138'
$ws.Range("E49").Value = '  -2.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '65.00'
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.21'
$ws.Range("E51").Value = '  -2.78%  '
